# Add a new "code line" column (E) and its supporting "file" column (F)
# to the header row, and fill in the code-line value for the
# "Priority assessment" row (row 2), per commit:
# "added code line of priority assesment"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "code line"
$ws.Range("F1").Value = "file"

$ws.Range("E2").Value = 343
$ws.Range("F2").Value = "value_functions.nls"
